$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative OOXML diff.
# Columns D (Price) and E (Volume(1h)) hold numeric-looking / percent-looking
# text that Excel would otherwise auto-convert to a number on assignment, so
# we force a text format, assign the literal string, then restore the default
# "Normal" style so the saved XML keeps matching the workbook's original look.
function Set-TextValue($rangeRef, $value) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '307.87'
Set-TextValue "E2" '1.42%'
Set-TextValue "D3" '39.39'
Set-TextValue "E3" '10.32%'
Set-TextValue "D4" '5.099'
Set-TextValue "E4" '1.35%'
Set-TextValue "D5" '0.08149'
Set-TextValue "E5" '3.23%'
Set-TextValue "D6" '1.990'
Set-TextValue "E6" '7.99%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue "D7" '7.901'
Set-TextValue "E7" '1.46%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D8" '0.9301'
Set-TextValue "E8" '1.13%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D9" '0.1416'
Set-TextValue "E9" '5.50%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D10" '0.1942'
Set-TextValue "E10" '2.38%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D11" '0.09255'
Set-TextValue "E11" '1.62%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D12" '0.03503'
Set-TextValue "E12" '1.14%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D13" '0.09821'
Set-TextValue "E13" '-0.09%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D14" '0.001409'
Set-TextValue "E14" '0.36%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D15" '0.005963'
Set-TextValue "E15" '-2.96%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D16" '3.947'
Set-TextValue "E16" '5.91%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D17" '4.178'
Set-TextValue "E17" '1.78%'
Set-TextValue "E18" '2.67%'
Set-TextValue "D19" '0.3452'
Set-TextValue "E19" '0.39%'
Set-TextValue "D20" '0.1303'
Set-TextValue "E20" '-0.46%'
Set-TextValue "D21" '4.809'
Set-TextValue "E21" '-6.93%'
Set-TextValue "E22" '19.47%'
Set-TextValue "D23" '0.04474'
Set-TextValue "E23" '1.58%'
Set-TextValue "D24" '0.001245'
Set-TextValue "E24" '0.86%'
Set-TextValue "E25" '-9.65%'
Set-TextValue "E27" '-0.04%'
Set-TextValue "D39" '0.02109'
Set-TextValue "E39" '8.59%'
Set-TextValue "D40" '0.05162'
Set-TextValue "E40" '1.62%'
Set-TextValue "D41" '0.007462'
Set-TextValue "E41" '-1.72%'
Set-TextValue "D42" '0.01018'
Set-TextValue "E42" '0.01%'
Set-TextValue "E43" '1.68%'
Set-TextValue "E44" '-0.97%'
Set-TextValue "D45" '0.009669'
Set-TextValue "E45" '-5.08%'
Set-TextValue "D46" '0.00006326'
Set-TextValue "E46" '2.70%'
Set-TextValue "E47" '0.00%'
Set-TextValue "E48" '-0.64%'
Set-TextValue "E49" '-3.46%'
Set-TextValue "D50" '0.00002103'
Set-TextValue "E50" '0.00%'
Set-TextValue "E51" '0.00%'
